$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.388527154922485
$ws.Range("B1").Value = 2.551299571990967
$ws.Range("C1").Value = 1.954998970031738
$ws.Range("D1").Value = 1.751769304275513
$ws.Range("E1").Value = 1.57250702381134
